# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary / stats block (rows 10-12) ---
# Row headers (No. / Marking / Total) get the "mtitleStyle" formatting (s="4"),
# matching the rest of the header column.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Right / Wrong / Not-Attempt / Max counters
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# Marking scheme - correct C11 which was incorrectly stored as text "-1"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Totals
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "78/112"

# --- Remove the third (unused) Student Ans / Correct Ans block in columns G:H ---
$ws.Columns("G:H").Delete()

# --- Second answer block (columns D:E) only keeps the first 3 questions (rows 16-18) ---
# Clear out rows 19-40 in D:E entirely so the cells (and the sheet dimension) shrink.
$ws.Range("D19:E40").Clear()

# Student answers for columns D (rows 16-18), paired with the already-present correct
# answers in column E, colored green (correctStyle) when matching, red (incorrectStyle)
# otherwise.
$colD = @{
    16 = "Option A"
    17 = "Option C"
    18 = "Option B"
}
foreach ($r in $colD.Keys) {
    $studentAns = $colD[$r]
    $correctAns = [string]$ws.Range("E$r").Value2
    $ws.Range("D$r").Value = $studentAns
    if ($studentAns -eq $correctAns) {
        $ws.Range("D$r").Style = "correctStyle"
    } else {
        $ws.Range("D$r").Style = "incorrectStyle"
    }
}

# --- First answer block (column A, rows 16-40) now holds the student's actual answers ---
$colA = @{
    16 = $null
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = $null
    23 = "Option D"
    24 = "Option A"
    25 = $null
    26 = "Option D"
    27 = "Option A"
    28 = $null
    29 = "Option D"
    30 = "Option B"
    31 = "Option D"
    32 = "Option C"
    33 = "Option D"
    34 = $null
    35 = "Option D"
    36 = $null
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}
foreach ($r in $colA.Keys) {
    $studentAns = $colA[$r]
    $correctAns = [string]$ws.Range("B$r").Value2
    if ($null -eq $studentAns) {
        # Not attempted - leave blank, keep the default "normalStyle" look
        $ws.Range("A$r").Value = ""
        $ws.Range("A$r").Style = "normalStyle"
    } else {
        $ws.Range("A$r").Value = $studentAns
        if ($studentAns -eq $correctAns) {
            $ws.Range("A$r").Style = "correctStyle"
        } else {
            $ws.Range("A$r").Style = "incorrectStyle"
        }
    }
}
